$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.644
$ws.Range("B2").Value = 0.044
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0.531
